# Apply the edits described by the commit diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the QPU parameter inputs (row 8) that drive the cost formulas below.
$ws.Range("B8").Value = 220
$ws.Range("C8").Value = 120
$ws.Range("D8").Value = 300

# Move the active selection from F8 to F10, matching the saved cursor position.
[void]$ws.Range("F10").Select()
